$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.436.17'
$ws.Range("E2").Value = '  -1.13%  '

$ws.Range("D3").Value = '2.686.17'
$ws.Range("E3").Value = '  -2.80%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -2.79%  '

$ws.Range("E9").Value = '  -3.86%  '

$ws.Range("E10").Value = '  -1.05%  '

$ws.Range("E11").Value = '  -4.30%  '

$ws.Range("E12").Value = '  -9.30%  '

$ws.Range("D13").Value = '3.162.13'
$ws.Range("E13").Value = '  -2.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("D15").Value = '63.296.79'
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("E16").Value = '  -3.88%  '

$ws.Range("D17").Value = '2.687.13'

$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("E19").Value = '  -5.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("E23").Value = '  -3.57%  '

$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("E25").Value = '  -1.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.19%  '

$ws.Range("D28").Value = '0.0₃0860'
$ws.Range("E28").Value = '  -5.31%  '

$ws.Range("E29").Value = '  +0.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.02%  '

$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("E34").Value = '  -2.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.08%  '

$ws.Range("E36").Value = '  -4.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '343.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.942'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.16%  '

$ws.Range("E40").Value = '  -2.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.14'
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = '  -6.10%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.35%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.60%  '

$ws.Range("E45").Value = '  -1.75%  '

$ws.Range("E46").Value = '  -4.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0973'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.67%  '

$ws.Range("E51").Value = '  -4.59%  '
